# Applies the cryptocurrency price/volume refresh described in the commit:
# "Updated cryptos list on Fri Aug  9 18:33:39 UTC 2024 with GitHub Actions"
#
# For each changed cell, the new value is written back. Cells in column D
# whose new text would otherwise be re-interpreted by Excel as a plain
# number (losing formatting such as trailing zeros, e.g. "1.00" -> 1) are
# first forced to Text number format so the literal string is preserved,
# exactly like the original inline strings in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(

    @{ Cell = "D2"; Value = '60.032.36'; ForceText = $false }
    @{ Cell = "E2"; Value = '  +1.61%  '; ForceText = $false }
    @{ Cell = "D3"; Value = '2.583.68'; ForceText = $false }
    @{ Cell = "E3"; Value = '  -0.07%  '; ForceText = $false }
    @{ Cell = "E4"; Value = '  -0.22%  '; ForceText = $false }
    @{ Cell = "D5"; Value = '505.70'; ForceText = $true }
    @{ Cell = "E5"; Value = '  +0.22%  '; ForceText = $false }
    @{ Cell = "D6"; Value = '152.57'; ForceText = $true }
    @{ Cell = "E6"; Value = '  -2.93%  '; ForceText = $false }
    @{ Cell = "D7"; Value = '0.998'; ForceText = $true }
    @{ Cell = "E7"; Value = '  +0.14%  '; ForceText = $false }
    @{ Cell = "D8"; Value = '0.576'; ForceText = $true }
    @{ Cell = "E8"; Value = '  -7.60%  '; ForceText = $false }
    @{ Cell = "D9"; Value = '2.588.64'; ForceText = $false }
    @{ Cell = "E9"; Value = '  +0.12%  '; ForceText = $false }
    @{ Cell = "D10"; Value = '6.62'; ForceText = $true }
    @{ Cell = "E10"; Value = '  +6.89%  '; ForceText = $false }
    @{ Cell = "E11"; Value = '  -0.09%  '; ForceText = $false }
    @{ Cell = "E12"; Value = '  +1.55%  '; ForceText = $false }
    @{ Cell = "D13"; Value = '0.127'; ForceText = $true }
    @{ Cell = "E13"; Value = '  +0.78%  '; ForceText = $false }
    @{ Cell = "D14"; Value = '3.037.09'; ForceText = $false }
    @{ Cell = "E14"; Value = '  +0.17%  '; ForceText = $false }
    @{ Cell = "D15"; Value = '60.069.60'; ForceText = $false }
    @{ Cell = "E15"; Value = '  +1.70%  '; ForceText = $false }
    @{ Cell = "D16"; Value = '21.47'; ForceText = $true }
    @{ Cell = "E16"; Value = '  -1.34%  '; ForceText = $false }
    @{ Cell = "D17"; Value = '0.0000140'; ForceText = $true }
    @{ Cell = "E17"; Value = '  +1.74%  '; ForceText = $false }
    @{ Cell = "D18"; Value = '2.581.37'; ForceText = $false }
    @{ Cell = "E18"; Value = '  -0.21%  '; ForceText = $false }
    @{ Cell = "D19"; Value = '4.83'; ForceText = $true }
    @{ Cell = "E19"; Value = '  +1.89%  '; ForceText = $false }
    @{ Cell = "D20"; Value = '344.89'; ForceText = $true }
    @{ Cell = "E20"; Value = '  +3.00%  '; ForceText = $false }
    @{ Cell = "D21"; Value = '10.38'; ForceText = $true }
    @{ Cell = "E21"; Value = '  +0.44%  '; ForceText = $false }
    @{ Cell = "D22"; Value = '6.14'; ForceText = $true }
    @{ Cell = "E22"; Value = '  +1.40%  '; ForceText = $false }
    @{ Cell = "E23"; Value = '  -0.80%  '; ForceText = $false }
    @{ Cell = "D24"; Value = '60.00'; ForceText = $true }
    @{ Cell = "E24"; Value = '  +0.25%  '; ForceText = $false }
    @{ Cell = "D25"; Value = '0.421'; ForceText = $true }
    @{ Cell = "E25"; Value = '  +1.35%  '; ForceText = $false }
    @{ Cell = "D26"; Value = '0.165'; ForceText = $true }
    @{ Cell = "E26"; Value = '  -0.73%  '; ForceText = $false }
    @{ Cell = "B27"; Value = 'Binance-PegBSC-USD'; ForceText = $false }
    @{ Cell = "C27"; Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; ForceText = $false }
    @{ Cell = "D27"; Value = '1.00'; ForceText = $true }
    @{ Cell = "E27"; Value = '  -0.15%  '; ForceText = $false }
    @{ Cell = "B28"; Value = 'PEPE'; ForceText = $false }
    @{ Cell = "C28"; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; ForceText = $false }
    @{ Cell = "D28"; Value = '0.0₃0842'; ForceText = $false }
    @{ Cell = "E28"; Value = '  +2.06%  '; ForceText = $false }
    @{ Cell = "B29"; Value = 'InternetComputer(DFINITY)'; ForceText = $false }
    @{ Cell = "C29"; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText = $false }
    @{ Cell = "D29"; Value = '7.35'; ForceText = $true }
    @{ Cell = "E29"; Value = '  -0.67%  '; ForceText = $false }
    @{ Cell = "B30"; Value = 'USDe'; ForceText = $false }
    @{ Cell = "C30"; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'; ForceText = $false }
    @{ Cell = "D30"; Value = '1.00'; ForceText = $true }
    @{ Cell = "E30"; Value = '  -0.03%  '; ForceText = $false }
    @{ Cell = "B31"; Value = 'EthereumClassic'; ForceText = $false }
    @{ Cell = "C31"; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; ForceText = $false }
    @{ Cell = "D31"; Value = '19.32'; ForceText = $true }
    @{ Cell = "E31"; Value = '  -0.15%  '; ForceText = $false }
    @{ Cell = "B32"; Value = 'Monero'; ForceText = $false }
    @{ Cell = "C32"; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; ForceText = $false }
    @{ Cell = "D32"; Value = '152.84'; ForceText = $true }
    @{ Cell = "E32"; Value = '  -2.97%  '; ForceText = $false }
    @{ Cell = "B33"; Value = 'PancakeSwap'; ForceText = $false }
    @{ Cell = "C33"; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; ForceText = $false }
    @{ Cell = "D33"; Value = '1.55'; ForceText = $true }
    @{ Cell = "E33"; Value = '  -1.00%  '; ForceText = $false }
    @{ Cell = "B34"; Value = 'Aptos'; ForceText = $false }
    @{ Cell = "C34"; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; ForceText = $false }
    @{ Cell = "D34"; Value = '5.76'; ForceText = $true }
    @{ Cell = "E34"; Value = '  +5.36%  '; ForceText = $false }
    @{ Cell = "B35"; Value = 'NEARProtocol'; ForceText = $false }
    @{ Cell = "C35"; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; ForceText = $false }
    @{ Cell = "D35"; Value = '3.98'; ForceText = $true }
    @{ Cell = "E35"; Value = '  +1.32%  '; ForceText = $false }
    @{ Cell = "B36"; Value = 'ImmutableX'; ForceText = $false }
    @{ Cell = "C36"; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; ForceText = $false }
    @{ Cell = "D36"; Value = '1.18'; ForceText = $true }
    @{ Cell = "E36"; Value = '  -0.45%  '; ForceText = $false }
    @{ Cell = "B37"; Value = 'SuiNetwork'; ForceText = $false }
    @{ Cell = "C37"; Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'; ForceText = $false }
    @{ Cell = "D37"; Value = '0.851'; ForceText = $true }
    @{ Cell = "E37"; Value = '  +16.91%  '; ForceText = $false }
    @{ Cell = "B38"; Value = 'Fetch.AI'; ForceText = $false }
    @{ Cell = "C38"; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; ForceText = $false }
    @{ Cell = "D38"; Value = '0.846'; ForceText = $true }
    @{ Cell = "E38"; Value = '  -0.75%  '; ForceText = $false }
    @{ Cell = "B39"; Value = 'Filecoin'; ForceText = $false }
    @{ Cell = "C39"; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText = $false }
    @{ Cell = "D39"; Value = '3.75'; ForceText = $true }
    @{ Cell = "E39"; Value = '  -0.16%  '; ForceText = $false }
    @{ Cell = "E40"; Value = '  +1.73%  '; ForceText = $false }
    @{ Cell = "B41"; Value = 'OKB'; ForceText = $false }
    @{ Cell = "C41"; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; ForceText = $false }
    @{ Cell = "D41"; Value = '35.84'; ForceText = $true }
    @{ Cell = "E41"; Value = '  +2.36%  '; ForceText = $false }
    @{ Cell = "B42"; Value = 'Bittensor'; ForceText = $false }
    @{ Cell = "C42"; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; ForceText = $false }
    @{ Cell = "D42"; Value = '292.67'; ForceText = $true }
    @{ Cell = "E42"; Value = '  +1.09%  '; ForceText = $false }
    @{ Cell = "B43"; Value = 'Stellar'; ForceText = $false }
    @{ Cell = "C43"; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; ForceText = $false }
    @{ Cell = "D43"; Value = '0.0996'; ForceText = $true }
    @{ Cell = "E43"; Value = '  -1.78%  '; ForceText = $false }
    @{ Cell = "B44"; Value = 'Mantle'; ForceText = $false }
    @{ Cell = "C44"; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; ForceText = $false }
    @{ Cell = "D44"; Value = '0.615'; ForceText = $true }
    @{ Cell = "E44"; Value = '  -1.69%  '; ForceText = $false }
    @{ Cell = "B45"; Value = 'Hedera'; ForceText = $false }
    @{ Cell = "C45"; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; ForceText = $false }
    @{ Cell = "D45"; Value = '0.0557'; ForceText = $true }
    @{ Cell = "E45"; Value = '  -0.82%  '; ForceText = $false }
    @{ Cell = "D46"; Value = '0.997'; ForceText = $true }
    @{ Cell = "E46"; Value = '  +0.26%  '; ForceText = $false }
    @{ Cell = "B47"; Value = 'EnergySwap'; ForceText = $false }
    @{ Cell = "C47"; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false }
    @{ Cell = "D47"; Value = '19.77'; ForceText = $true }
    @{ Cell = "E47"; Value = '  +2.13%  '; ForceText = $false }
    @{ Cell = "B48"; Value = 'RenderToken'; ForceText = $false }
    @{ Cell = "C48"; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText = $false }
    @{ Cell = "D48"; Value = '4.86'; ForceText = $true }
    @{ Cell = "E48"; Value = '  +1.33%  '; ForceText = $false }
    @{ Cell = "B49"; Value = 'VeChain'; ForceText = $false }
    @{ Cell = "C49"; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $false }
    @{ Cell = "D49"; Value = '0.0233'; ForceText = $true }
    @{ Cell = "E49"; Value = '  -1.03%  '; ForceText = $false }
    @{ Cell = "B50"; Value = 'WhiteBITCoin'; ForceText = $false }
    @{ Cell = "C50"; Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'; ForceText = $false }
    @{ Cell = "D50"; Value = '10.29'; ForceText = $true }
    @{ Cell = "E50"; Value = '  -0.19%  '; ForceText = $false }
    @{ Cell = "B51"; Value = 'Maker'; ForceText = $false }
    @{ Cell = "C51"; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; ForceText = $false }
    @{ Cell = "D51"; Value = '2.004.25'; ForceText = $false }
    @{ Cell = "E51"; Value = '  +0.99%  '; ForceText = $false }
)

foreach ($change in $changes) {
    $range = $ws.Range($change.Cell)
    if ($change.ForceText) {
        $range.NumberFormat = "@"
    }
    $range.Value = $change.Value
}
